# Update stats for 2025-11 (row 24)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = 6352
$ws.Range("C24").Value = 999
$ws.Range("D24").Value = 5955330
$ws.Range("E24").Value = 937.5519521410579
$ws.Range("F24").Value = 8.285032390044321
$ws.Range("G24").Value = 3.523316062176174
$ws.Range("H24").Value = 26.15695725009597
